$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Turn the lone "hola" paragraph into the memo header line, then
#    append the rest of the memo body as new paragraphs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("hola", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "titulo: lliii", 2)

$lines = @(
    "REF: et/SR",
    "Contenido:",
    "sdf",
    " ",
    "f._______________________________",
    "Samuel Rabanales",
    "Jefe de turno",
    "Subcomisaria 41-31 San Juan Ostuncalco"
)

foreach ($line in $lines) {
    $d.Paragraphs.Last.Range.InsertParagraphAfter()
    $d.Paragraphs.Last.Range.InsertBefore($line)
}

# ------------------------------------------------------------------
# 2. Resize the page from the default A4-ish size to a US Letter /
#    legal-tall 12240 x 18720 twips (612 x 936 pt) canvas.
# ------------------------------------------------------------------
$ps = $d.PageSetup
$ps.PageWidth = 612
$ps.PageHeight = 936

# ------------------------------------------------------------------
# 3. Update the document's default run formatting: lower-cased
#    "arial" font family (all script slots) at 12pt / 24 half-points
#    instead of Arial 10pt.
# ------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$f = $normalStyle.Font
$f.Name = "arial"
$f.NameFarEast = "arial"
$f.NameBi = "arial"
$f.Size = 12
$f.SizeBi = 12
